$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally held two "Sending cluster" blocks (FAPs, then
# Resolving-Mac) x four "Target cluster" rows each. The new TPM re-run
# drops the Resolving-Mac sending-cluster block entirely (rows 6-9) and
# refreshes the computed columns (I, J, M-T) for the remaining FAPs block.
$ws.Rows.Item(9).EntireRow.Delete() | Out-Null
$ws.Rows.Item(8).EntireRow.Delete() | Out-Null
$ws.Rows.Item(7).EntireRow.Delete() | Out-Null
$ws.Rows.Item(6).EntireRow.Delete() | Out-Null

# Refreshed "specificity of expression" pair (columns I, J) is now 1 for
# every remaining row (previously 0.4255051361566043).
$ws.Range("I2:J5").Value = 1

# Row 2: FAPs -> ECs
$ws.Range("M2").Value = 14.440165
$ws.Range("N2").Value = 43.320495
$ws.Range("O2").Value = 0.1441015470002482
$ws.Range("P2").Value = 0.1441015470002482
$ws.Range("Q2").Value = 1.690765226131667
$ws.Range("R2").Value = 15.216887035185
$ws.Range("S2").Value = 0.1441015470002482
$ws.Range("T2").Value = 0.1441015470002482

# Row 3: FAPs -> FAPs
$ws.Range("O3").Value = 0.3846359116098663
$ws.Range("P3").Value = 0.3846359116098662
$ws.Range("S3").Value = 0.3846359116098663
$ws.Range("T3").Value = 0.3846359116098662

# Row 4: FAPs -> MuSCs
$ws.Range("M4").Value = 21.954262
$ws.Range("N4").Value = 65.862786
$ws.Range("O4").Value = 0.2190863551385157
$ws.Range("P4").Value = 0.2190863551385156
$ws.Range("Q4").Value = 2.570573310968666
$ws.Range("R4").Value = 23.135159798718
$ws.Range("S4").Value = 0.2190863551385157
$ws.Range("T4").Value = 0.2190863551385156

# Row 5: FAPs -> Resolving-Mac
$ws.Range("M5").Value = 25.27013633333333
$ws.Range("N5").Value = 75.81040899999999
$ws.Range("O5").Value = 0.2521761862513699
$ws.Range("P5").Value = 0.2521761862513699
$ws.Range("Q5").Value = 2.958821299618555
$ws.Range("R5").Value = 26.629391696567
$ws.Range("S5").Value = 0.2521761862513699
$ws.Range("T5").Value = 0.2521761862513699
